# "updated entire analysis with newest numbers from 7th October 2024"
# Append 10 new rows (55-64) of 2024 churn data to Sheet1.
#
# Column layout: A=CVR B=Year C=Beløb(TCV) D=Løsning E=Opsagt dato F=Årsag
#                G=Ny leverandør H=Quarter I=TCV_range
#
# New values are written column-by-column (all of column A for the new
# rows, then D, then G, then H) so that newly-introduced shared strings
# land in the workbook's shared-string table in the same order they were
# first encountered when the source data was pasted in (CVR numbers
# first, then the new "Løsning" value, then the new "Ny leverandør"
# value, then the new quarters) rather than row-by-row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cvr    = @("25227832","46481410","71133214","31943140","35388451","19691403","29414815","15504749","39427168","12316100")
$losning = @("Visma Løn","Visma Time","Visma Løn","Kundeforhold","Visma Løn","SKL SE","EasyCruit","Visma Løn","EasyCruit","Visma Løn")
$quarter = @("2024Q1","2024Q2","2024Q2","2023Q4","2024Q2","2024Q2","2024Q2","2024Q2","2024Q3","2024Q3")
$year    = @(2024,2024,2024,2024,2024,2024,2024,2024,2024,2024)
$belob   = @(32083.08,33984,25747.01,31550,29760,21012,37272,20388,31956,24988)
$leverandor = @{ 3 = "DataLøn"; 5 = "DataLøn"; 9 = "Workday"; 10 = "DataLøn" }
$year2024 = @(2024,6,3,2024,5,8,2023,11,2,2024,5,17,2024,6,5,2024,6,10,2024,6,13,2024,8,30,2024,9,24)
$dates = @(
    @(2024,3,8),
    @(2024,6,3),
    @(2024,5,8),
    @(2023,11,2),
    @(2024,5,17),
    @(2024,6,5),
    @(2024,6,10),
    @(2024,6,13),
    @(2024,8,30),
    @(2024,9,24)
)
$range = "20000-40000"

$startRow = 55
$n = 10

# --- Column A: CVR numbers (stored as text, matching the existing column).
# The values are digit-only, so a plain Value assignment would be coerced
# to a number by Excel; format the cell as Text first, assign, then put
# the cell style back to Normal so the saved cell carries no explicit
# style (matching the rest of the column) while the stored value stays a
# string.
for ($i = 0; $i -lt $n; $i++) {
    $r = $startRow + $i
    $cell = $ws.Cells.Item($r, 1)
    $cell.NumberFormat = "@"
    $cell.Value = $cvr[$i]
    $cell.Style = "Normal"
}

# --- Column D: Løsning
for ($i = 0; $i -lt $n; $i++) {
    $r = $startRow + $i
    $ws.Cells.Item($r, 4).Value = $losning[$i]
}

# --- Column G: Ny leverandør (only set on the rows that have one)
foreach ($i in $leverandor.Keys) {
    $r = $startRow + $i - 1
    $ws.Cells.Item($r, 7).Value = $leverandor[$i]
}

# --- Column H: Quarter
for ($i = 0; $i -lt $n; $i++) {
    $r = $startRow + $i
    $ws.Cells.Item($r, 8).Value = $quarter[$i]
}

# --- Remaining non-string columns: B (Year), C (Beløb), E (Opsagt dato), I (TCV_range)
for ($i = 0; $i -lt $n; $i++) {
    $r = $startRow + $i
    $ws.Cells.Item($r, 2).Value = $year[$i]
    $ws.Cells.Item($r, 3).Value = $belob[$i]

    $d = $dates[$i]
    $eCell = $ws.Cells.Item($r, 5)
    $eCell.NumberFormat = "YYYY-MM-DD HH:MM:SS"
    $eCell.Value = (Get-Date -Year $d[0] -Month $d[1] -Day $d[2] -Hour 0 -Minute 0 -Second 0)

    $ws.Cells.Item($r, 9).Value = $range
}
